$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.434.56'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '3.048.79'
$ws.Range("E3").Value = '  +4.41%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.76'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.68'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("E9").Value = '  +6.20%  '
$ws.Range("D10").Value = '3.047.65'
$ws.Range("E10").Value = '  +4.34%  '
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.16'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.17%  '
$ws.Range("D14").Value = '3.609.58'
$ws.Range("E14").Value = '  +4.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.36'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.60%  '
$ws.Range("D16").Value = '76.381.13'
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("E17").Value = '  +2.66%  '
$ws.Range("D18").Value = '3.054.03'
$ws.Range("E18").Value = '  +4.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.57'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.88%  '
$ws.Range("E20").Value = '  +4.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.52'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("E23").Value = '  +2.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.66'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.11%  '
$ws.Range("D25").Value = '3.207.63'
$ws.Range("E25").Value = '  +4.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.43'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.94'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("E29").Value = '  +3.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.32'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +7.65%  '
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '506.66'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("E34").Value = '  +6.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.84'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.12'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.387'
$ws.Range("D38").ClearFormats()
$ws.Range("E39").Value = '  +2.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '191.29'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.65%  '
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.20'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.08%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.797'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +21.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.28'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +7.43%  '
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("E49").Value = '  +5.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.611'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +6.90%  '
$ws.Range("E51").Value = '  +5.61%  '
